# Helper: set a cell's value while preventing Excel from auto-converting
# numeric-looking text (e.g. "214.68") into a real number. We prefix with
# an apostrophe to force text, then reset the style back to "Normal" so no
# visible/structural style change is left behind on the cell.
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub3 = [string][char]8323
$sub6 = [string][char]8326

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "27.838.99"
$ws.Range("E2").Value = "  +2.57%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.663.60"
$ws.Range("E3").Value = "  -0.75%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "214.68"
$ws.Range("E5").Value = "  +0.26%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.48%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "23.45"
$ws.Range("E8").Value = "  +2.60%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.260"
$ws.Range("E9").Value = "  -0.28%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0620"
$ws.Range("E10").Value = "  -0.16%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.0878"
$ws.Range("E11").Value = "  -1.32%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "1.900.40"
$ws.Range("E12").Value = "  -0.70%  "

# Row 13 - WrappedEther
Set-TextValue $ws.Range("D13") "1.666.94"
$ws.Range("E13").Value = "  -0.58%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "4.14"
$ws.Range("E14").Value = "  -1.64%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("D15") "0.548"

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "65.97"
$ws.Range("E16").Value = "  -0.82%  "

# Row 17 - BitcoinCash
Set-TextValue $ws.Range("D17") "248.99"
$ws.Range("E17").Value = "  +5.88%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "27.804.92"
$ws.Range("E18").Value = "  +2.52%  "

# Row 19 - ShibaInu
$d19val = "0.0" + $sub3 + "0730"
Set-TextValue $ws.Range("D19") $d19val
$ws.Range("E19").Value = "  -1.54%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "7.49"
$ws.Range("E20").Value = "  -4.29%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.07%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "4.47"
$ws.Range("E22").Value = "  -1.57%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  -2.29%  "

# Row 24 - Toncoin
Set-TextValue $ws.Range("D24") "2.05"
$ws.Range("E24").Value = "  -1.73%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "146.64"
$ws.Range("E25").Value = "  -0.82%  "

# Row 26 - Cosmos
Set-TextValue $ws.Range("D26") "7.19"
$ws.Range("E26").Value = "  -3.77%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "16.26"
$ws.Range("E27").Value = "  -1.11%  "

# Rows 28 and 29 swap places: Stellar <-> BinanceUSD
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D29") "0.112"
$ws.Range("E29").Value = "  -0.61%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +6.11%  "

# Row 31 - Hedera
Set-TextValue $ws.Range("D31") "0.0499"
$ws.Range("E31").Value = "  -0.02%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "3.34"
$ws.Range("E32").Value = "  -0.77%  "

# Row 33 - Maker
Set-TextValue $ws.Range("D33") "1.432.70"
$ws.Range("E33").Value = "  -7.34%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D34") "3.13"
$ws.Range("E34").Value = "  -3.13%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -6.54%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.03%  "

# Row 37 - ARBITRUM
Set-TextValue $ws.Range("D37") "0.927"
$ws.Range("E37").Value = "  -2.09%  "

# Row 38 - ImmutableX
Set-TextValue $ws.Range("D38") "0.582"
$ws.Range("E38").Value = "  -4.22%  "

# Row 39 - VeChain
Set-TextValue $ws.Range("D39") "0.0169"
$ws.Range("E39").Value = "  -1.65%  "

# Row 40 - WEMIXToken
Set-TextValue $ws.Range("D40") "1.05"
$ws.Range("E40").Value = "  -2.14%  "

# Row 41 - Aave
Set-TextValue $ws.Range("D41") "69.71"
$ws.Range("E41").Value = "  +0.00%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.03%  "

# Rows 43 and 44 swap places: MXToken <-> FraxShare
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D43") "5.41"
$ws.Range("E43").Value = "  -6.45%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D44") "2.22"
$ws.Range("E44").Value = "  -1.51%  "

# Row 45 - RocketPoolETH
Set-TextValue $ws.Range("D45") "1.806.50"
$ws.Range("E45").Value = "  -1.14%  "

# Row 46 - TrustWalletToken
$ws.Range("E46").Value = "  +0.73%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +3.74%  "

# Row 48 - Quant
Set-TextValue $ws.Range("D48") "89.15"
$ws.Range("E48").Value = "  -0.60%  "

# Row 49 - BabyDogeCoin
$d49val = "0.0" + $sub6 + "0110"
Set-TextValue $ws.Range("D49") $d49val
$ws.Range("E49").Value = "  -1.06%  "

# Row 50 - Algorand
Set-TextValue $ws.Range("D50") "0.101"
$ws.Range("E50").Value = "  -2.84%  "

# Row 51 - EnergySwap
Set-TextValue $ws.Range("D51") "7.81"
$ws.Range("E51").Value = "  -5.35%  "
